# Commit: "search for test no need to install sqllite on Dockerfile"
# The Testdaten sheet held a sample search term "Testanalyst" in A2;
# update it to the shorter search term "Test" and leave the selection
# cursor resting on A5 (matches the saved sheet view in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdaten")

$ws.Range("A2").Value = "Test"

$ws.Range("A5").Select()
